$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: "... reason is lack of motivation, which is ..." -> remove "is "
# ---------------------------------------------------------------------------
$r1 = $d.Content
$old1 = " reason is lack of motivation, which is "
$new1 = " reason is lack of motivation, which "
$r1.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# ---------------------------------------------------------------------------
# Edit 2/3: shift "they will " across the run boundary (net text unchanged)
# and drop the stale lastRenderedPageBreak hint that no longer lines up with
# the page after the text above was shortened.
# ---------------------------------------------------------------------------
$r2 = $d.Content
$old2 = "as they do not have all the information available to them, they will find it hard to be motivated and do their own research"
$r2.Find.Execute($old2, $false, $false, $false, $false, $false, $true, 1, $false, $old2, 2) | Out-Null

# ---------------------------------------------------------------------------
# New paragraph: "Aims and Objectives" heading (bold) inserted right before
# the "This project aims to ..." paragraph.
# ---------------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute("This project aims to") | Out-Null
$targetPara = $r3.Paragraphs.First
$targetPara.Range.InsertParagraphBefore()

# Re-locate the (now shifted) "This project aims to" paragraph and grab the
# freshly inserted blank paragraph immediately preceding it.
$r3b = $d.Content
$r3b.Find.Execute("This project aims to") | Out-Null
$targetPara2 = $r3b.Paragraphs.First
$headingPara = $targetPara2.Previous

$hr = $headingPara.Range
$hr.InsertBefore("Aims and Objectives")
$hr.Font.Bold = 1
$hr.Font.BoldBi = 1

# ---------------------------------------------------------------------------
# Edit 4: " motivated over time, " -> " motivated and engaged over time, "
# ---------------------------------------------------------------------------
$r4 = $d.Content
$old4 = "to keep the users motivated over time, helping"
$new4 = "to keep the users motivated and engaged over time, helping"
$r4.Find.Execute($old4, $false, $false, $false, $false, $false, $true, 1, $false, $new4, 2) | Out-Null
